$wb = $excel.ActiveWorkbook

# 1. Fix typo "REGREESOR" -> "REGRESSOR" in the header cell of HU_E1a.
$wsE1a = $wb.Worksheets.Item("HU_E1a")
$wsE1a.Range("A1").Value = "REGRESSOR"

# 2. Widen column A of HU_E2a (to fit the longer "_Medium"/"_Low" labels);
#    this width carries over to the duplicate created below.
$wsE2a = $wb.Worksheets.Item("HU_E2a")
$wsE2a.Columns.Item(1).ColumnWidth = 25.736979166666668

# 3. Duplicate the HU_E2a sheet, placing the copy right after it -> "HU_E2a (2)".
$wsE2a.Copy($null, $wsE2a)

# Copying activates the new sheet; restore HU_E2a as the active/selected tab.
$wb.Worksheets.Item("HU_E2a").Activate()

# 4. Trim the original HU_E2a sheet down to the smaller A1:I8 block, removing
#    the "_Medium" rows/columns (the full data now lives on "HU_E2a (2)").
$wsE2a.Range("J1:P8").Clear()
$wsE2a.Range("A9:P15").Clear()
